$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WEEK1")

$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# ---- Row 10 ----
$ws.Range("A10").Value = 1
$ws.Range("A10").HorizontalAlignment = $xlCenter

$ws.Range("B10").Value = "Individual"
$ws.Range("B10").HorizontalAlignment = $xlCenter

$ws.Range("C10").Value = "Tutorial Part 7 searching, there are two ways to create search, one is to modify URL, another is to add table into view. But the later way is the best to help users operate.`nUnderstood the way of adding search function.`n"
$ws.Range("C10").WrapText = $true

$ws.Range("D10").Value = "29/7/2020"

$ws.Rows.Item(10).RowHeight = 62.4

# ---- Row 11 ----
$ws.Range("A11").Value = 1
$ws.Range("A11").HorizontalAlignment = $xlCenter

$ws.Range("B11").Value = "Individual"
$ws.Range("B11").HorizontalAlignment = $xlCenter

$ws.Range("C11").Value = "Tutorial Part 8, studying how to add a new filed, when there will be a new field, all others will be updated based on new version. However, because of a tiny problem, my app can not process well."

$ws.Range("D11").Value = "29/7/2020"

# ---- Row 12 ----
$ws.Range("A12").Value = 2
$ws.Range("A12").HorizontalAlignment = $xlCenter

$ws.Range("B12").Value = "Individual"
$ws.Range("B12").HorizontalAlignment = $xlCenter

$ws.Range("C12").Value = "Tutorial Part 9 +10 , how to create validation rules that would be forced any time when users create or edit a movie. Meanwhile, explictly understood the details of fucntions of detail and delete. "
$ws.Range("C12").WrapText = $true

$ws.Range("D12").Value = "29/7/2020"

$ws.Rows.Item(12).RowHeight = 31.2

# ---- Row 13 ----
$ws.Range("A13").Value = 1
$ws.Range("A13").HorizontalAlignment = $xlCenter

$ws.Range("B13").Value = "Individual"
$ws.Range("B13").HorizontalAlignment = $xlCenter

$ws.Range("C13").Value = "Review the whole tutorial of MVC, tested and played some of them parts, trying to understand how to build my own website."

$ws.Range("D13").Value = "29/7/2020"

# ---- Final selection, mirrors the recorded cursor position after edits ----
$ws.Range("B17").Select()
